# Logged Week 15 and simulated Week 16
# Update the "Road" (R) row target-depth totals on both the OFF and DEF
# sheets to reflect the newly logged/simulated week's cumulative numbers.

$wb = $excel.ActiveWorkbook

# --- OFF sheet: row 3 is the "R" (Road) totals row ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 407
$wsOff.Range("C3").Value = 302
$wsOff.Range("D3").Value = 95
$wsOff.Range("E3").Value = 49
$wsOff.Range("F3").Value = 4

# --- DEF sheet: row 3 is the "R" (Road) totals row ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 403
$wsDef.Range("C3").Value = 283
$wsDef.Range("D3").Value = 96
$wsDef.Range("E3").Value = 46
